$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02)
$CValues = @(1.02657074974452,1.027957283232099,1.028852916152946,1.029229075492366,1.029292213074184,1.02885794383873,1.027039658922418,1.023823508096779,1.021670915329391,1.020736724342809,1.020389401062325,1.020463917821952,1.020708021125229,1.020858378194402,1.021732869578241,1.022280846641503,1.022600269854209,1.02270915077278,1.022222074889021,1.020636147763931,1.019637139242118,1.020166912481828,1.022248631940323,1.024656430660462)
$DValues = @(1.029448169621014,1.030450420371354,1.031097286991349,1.031368836065361,1.031414407350773,1.031100916986562,1.029787229756944,1.027459523327344,1.025898907296865,1.025221009557202,1.024968882635891,1.025022979509859,1.025200175330427,1.025309308271951,1.025943851703214,1.026341308166977,1.026572931518377,1.026651874223858,1.026298686214701,1.025148004595196,1.024422640080926,1.024807349553865,1.026317945878036,1.028062830194737)
$EValues = @(1.026761413817863,1.027949554640845,1.02871738055158,1.029039943658292,1.029094089994467,1.028721691558155,1.02716315695721,1.024409134746481,1.022567705555591,1.021769006131069,1.021472126078099,1.021535817327952,1.021744470195449,1.021873000465981,1.02262068377562,1.02308932127679,1.023362539835182,1.023455678440491,1.023039054366355,1.021683032895092,1.020829245222318,1.021281969926988,1.023061768220789,1.025122051540598)
$IValues = @(1.029876236138195,1.030171234074784,1.030359620905389,1.030438222198111,1.030451384753209,1.030360673522288,1.029976449700202,1.029280224870727,1.028803107682342,1.028593417266479,1.028515062214884,1.028531890784733,1.028586949943402,1.028620811817568,1.028816958787177,1.028939166698063,1.029010149943419,1.029034302803138,1.02892608583922,1.028570749289733,1.028344634122517,1.028464758621441,1.028931997438964,1.029462495579617)
$JValues = @(1.031733452048337,1.032757553557702,1.033418279209637,1.033695588424513,1.033742123004544,1.033421986431238,1.032079955184253,1.029700123608938,1.028103242524463,1.027409271391313,1.0271511175214,1.02720650977135,1.027387940134569,1.027499674593086,1.028149245786212,1.028556028950514,1.028793056734351,1.028873836188765,1.028512410028199,1.027334524001272,1.026591727211807,1.026985709244336,1.028532120289565,1.030317169452591)
$KValues = @(1.032262051586634,1.033071927337036,1.033593713704261,1.033812535047879,1.033849244666725,1.033596639714961,1.032536221981081,1.030650215273988,1.029380988395041,1.028828538046974,1.028622899020731,1.028667029005494,1.028811548743326,1.028900534454757,1.029417591993735,1.029741158440389,1.029929612927563,1.029993824274531,1.029706471433395,1.028769003303035,1.028177064904007,1.028491102277542,1.029722145852951,1.031139875397768)
$LValues = @(1.029583116537974,1.030577809862337,1.031219902839855,1.031489473537567,1.031534714330572,1.031223506287347,1.029919598821595,1.027609992884881,1.026061992769479,1.025389680709211,1.025139646413813,1.025193293614158,1.025369019098807,1.025477248465743,1.026106569010091,1.026500781769153,1.026730525021819,1.02680882869909,1.026458506609293,1.025317280852826,1.024597964623779,1.024979458224973,1.026477609535238,1.028208521096283)
$NValues = @(1.033198631891732,1.034224187742701,1.034885851700849,1.035163554726617,1.03521015539109,1.034889564187129,1.03354562710184,1.031162415892718,1.029563267053938,1.028868310402179,1.02860978992416,1.028665260837459,1.028846948852604,1.028958841986868,1.029609335645595,1.030016696488649,1.030254060879127,1.030334955049631,1.029973015622461,1.028793456862271,1.028049605216165,1.02844414674837,1.029992753874659,1.031780338012245)

$startRow = 2
for ($i = 0; $i -lt $BValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $BValues[$i]
    $ws.Cells.Item($row, 3).Value = $CValues[$i]
    $ws.Cells.Item($row, 4).Value = $DValues[$i]
    $ws.Cells.Item($row, 5).Value = $EValues[$i]
    $ws.Cells.Item($row, 9).Value = $IValues[$i]
    $ws.Cells.Item($row, 10).Value = $JValues[$i]
    $ws.Cells.Item($row, 11).Value = $KValues[$i]
    $ws.Cells.Item($row, 12).Value = $LValues[$i]
    $ws.Cells.Item($row, 14).Value = $NValues[$i]
}